# Anton_Rigin_SYRCoSE_2019.pptx - text tweaks on slides 14, 15, 27.
# Each edit below targets the exact run whose text changed (full run span),
# so PowerPoint doesn't have to split/merge runs and the surrounding
# <a:rPr> run-formatting (sub/superscript "+" "*" "*+" markers etc.) stays intact.

$p = $ppt.ActivePresentation

function Get-ShapeById($Slide, $TargetId, $FallbackIndex) {
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $shp = $Slide.Shapes.Item($i)
        if ($shp.Id -eq $TargetId) {
            return $shp
        }
    }
    return $Slide.Shapes.Item($FallbackIndex)
}

function Replace-RunText($TextRange, $OldRunText, $NewRunText) {
    $full = $TextRange.Text
    $idx = $full.IndexOf($OldRunText)
    if ($idx -lt 0) {
        throw "Could not locate expected run text: [$OldRunText]"
    }
    $sub = $TextRange.Characters($idx + 1, $OldRunText.Length)
    $sub.Text = $NewRunText
}

# Slide 14: "... the indexing data structure (B-tree, B+-tree, B*-tree or B*+-tree) ..."
# -> drop the leading "B-tree, " so it reads "... (B+-tree, B*-tree or B*+-tree) ..."
$slide14 = $p.Slides.Item(14)
$shape14 = Get-ShapeById $slide14 60 3
Replace-RunText $shape14.TextFrame.TextRange `
    "To develop and implement an algorithm that would allow selecting the indexing data structure (B-tree, B" `
    "To develop and implement an algorithm that would allow selecting the indexing data structure (B"

# Slide 15: "Selects from the B-tree and its modifications (B+-tree, ..." -> "Selects from the B-tree modifications (B+-tree, ..."
$slide15 = $p.Slides.Item(15)
$shape15 = Get-ShapeById $slide15 60 3
Replace-RunText $shape15.TextFrame.TextRange `
    "Selects from the B-tree and its modifications (B" `
    "Selects from the B-tree modifications (B"

# Slide 27: append ", however it has greater memory usage" to the B*+-tree complexity statement.
$slide27 = $p.Slides.Item(27)
$shape27 = Get-ShapeById $slide27 60 3
Replace-RunText $shape27.TextFrame.TextRange `
    "-tree has smaller computational complexity of keys insertion and deletion than B-tree" `
    "-tree has smaller computational complexity of keys insertion and deletion than B-tree, however it has greater memory usage"
